$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.85
$ws.Range("I3").Value = 5
$ws.Range("L3").Value = 6
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 4.75
$ws.Range("AX3").Value = 34

# Row 5 updates
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10

# Row 6 updates
$ws.Range("Q6").Value = 1.73
$ws.Range("R6").Value = 2.08
